$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 25614
$ws.Range("A3").Value = 25614
$ws.Range("A4").Value = 25614
$ws.Range("A5").Value = 25614
$ws.Range("A6").Value = 25614
$ws.Range("A7").Value = 25614
$ws.Range("A8").Value = 25614
$ws.Range("A9").Value = 25614
$ws.Range("A10").Value = 25614
$ws.Range("A11").Value = 25614
$ws.Range("A12").Value = 39401
$ws.Range("A13").Value = 39493
$ws.Range("A14").Value = 39583
$ws.Range("A15").Value = 39675
$ws.Range("A16").Value = 39767
$ws.Range("A17").Value = 39859
$ws.Range("A18").Value = 39948
$ws.Range("A19").Value = 40040
$ws.Range("A20").Value = 40132
$ws.Range("A21").Value = 40224
$ws.Range("A22").Value = 40313
$ws.Range("A23").Value = 40405
$ws.Range("A24").Value = 40497
$ws.Range("A25").Value = 40589
$ws.Range("A26").Value = 40678
$ws.Range("A27").Value = 40770
$ws.Range("A28").Value = 40862
$ws.Range("A29").Value = 40954
$ws.Range("A30").Value = 41044
$ws.Range("A31").Value = 41136
$ws.Range("A32").Value = 41228
$ws.Range("A33").Value = 41320
$ws.Range("A34").Value = 41409
$ws.Range("A35").Value = 41501
$ws.Range("A36").Value = 41593
$ws.Range("A37").Value = 41685
$ws.Range("A38").Value = 41774
$ws.Range("A39").Value = 41866
$ws.Range("A40").Value = 41958
$ws.Range("A41").Value = 42050
$ws.Range("A42").Value = 42139
$ws.Range("A43").Value = 42231
$ws.Range("A44").Value = 42323
$ws.Range("A45").Value = 42415
$ws.Range("A46").Value = 42505
$ws.Range("A47").Value = 42597
$ws.Range("A48").Value = 42689
$ws.Range("A49").Value = 42781
$ws.Range("A50").Value = 42870
$ws.Range("A51").Value = 42962
$ws.Range("A52").Value = 43054
$ws.Range("A53").Value = 43146
$ws.Range("A54").Value = 43235
$ws.Range("A55").Value = 43327
$ws.Range("A56").Value = 43419
$ws.Range("A57").Value = 43511
$ws.Range("A58").Value = 43600
$ws.Range("A59").Value = 43692
$ws.Range("A60").Value = 43784
$ws.Range("A61").Value = 43876
$ws.Range("A62").Value = 43966
$ws.Range("A63").Value = 44058
$ws.Range("A64").Value = 44150
$ws.Range("A65").Value = 44242
$ws.Range("A66").Value = 44331
$ws.Range("A67").Value = 44423
$ws.Range("A68").Value = 44515
$ws.Range("A69").Value = 44607
$ws.Range("A70").Value = 44696
$ws.Range("A71").Value = 44788
$ws.Range("A72").Value = 44880
$ws.Range("A73").Value = 44972
